# Q4 Financial Review - classification validation & restructure future labs
# Updates the customer sample data used on the "Top Customer Accounts" and
# "Premium Tier Customers" slides, and turns on autofit shrink-to-fit on the
# Premium Tier Customers body placeholder (PowerPoint applied autofit once the
# new text needed more room).

$p = $ppt.ActivePresentation

# --- Slide 2: "Top Customer Accounts" ---------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Top Customer Accounts:`r`n`r`nLauren Smith`r`n  CC: 6011-3336-6513-2206 (Discover)`r`n  SSN: 605-50-4048`r`n  Bank: 212524838 - 35928397`r`n`r`nRobert Rodriguez`r`n  CC: 3747-147727-40130 (Amex)`r`n  SSN: 339-93-9227`r`n  Bank: 310308176 - 10080973`r`n`r`nLauren Smith`r`n  CC: 3705-465003-28226 (Amex)`r`n  SSN: 583-06-4554`r`n  Bank: 661028505 - 69775754`r`n`r`nAmanda Wilson`r`n  CC: 5124-4082-4524-2787 (Mastercard)`r`n  SSN: 475-47-0764`r`n  Bank: 708747422 - 98891132`r`n`r`nJessica Moore`r`n  CC: 4844-4614-3142-2578 (Visa)`r`n  SSN: 854-45-5056`r`n  Bank: 344512868 - 23012863`r`n`r`n"

# --- Slide 3: "Premium Tier Customers" ---------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$shp3.TextFrame.TextRange.Text = "Premium Tier Customers:`r`n`r`nLauren Smith - lauren.smith@contoso.com`r`nLoyalty: RET-176207-4 | Phone: (392) 491-9134`r`n`r`nRobert Rodriguez - robert.rodriguez@contoso.com`r`nLoyalty: RET-050201-7 | Phone: (998) 566-6021`r`n`r`nLauren Smith - lauren.smith@gmail.com`r`nLoyalty: RET-605516-5 | Phone: (809) 627-8980`r`n`r`nAmanda Wilson - amanda.wilson@contoso.com`r`nLoyalty: RET-667621-0 | Phone: (364) 682-3815`r`n`r`nJessica Moore - jessica.moore@contoso.com`r`nLoyalty: RET-088512-3 | Phone: (730) 639-5655`r`n`r`nLauren Davis - lauren.davis@hotmail.com`r`nLoyalty: RET-561184-3 | Phone: (659) 762-7463`r`n`r`n"

# PowerPoint recalculated autofit for the larger block of text, shrinking the
# text to keep it inside the placeholder.
$shp3.TextFrame.AutoSize = 1
$shp3.TextFrame2.TextRange.ParagraphFormat.LineSpacingReduction = 20
$shp3.TextFrame.TextRange.Font.Size = $shp3.TextFrame.TextRange.Font.Size * 0.25
